# Update the Global Glider (GA05MOAS-GL002) Cal and Ingest sheet:
#   - CC_scattering_angle (row 2, col F) changes from 117 to 140
#   - CC_angular_resolution (row 4, col F) changes from 1.08 to 1.13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

$ws.Range("F2").Value = 140
$ws.Range("F4").Value = 1.13

# Leave the sheet with the same active selection seen in the authored workbook.
$ws.Range("F14").Select() | Out-Null
